$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values to reflect repulled data / recalculated means
$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -1
$ws.Range("F7").Value = -1
